# Minor fix for date: shift the Start (J) and End (K) date/time values
# for rows 2 through 28 forward by exactly 365 days (one year), leaving
# the time-of-day portion unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    foreach ($col in @("J", "K")) {
        $cell = $ws.Range("$col$row")
        $cell.Value2 = $cell.Value2 + 365
    }
}
